$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel
# are forced to Text format first, matching the original inline-string (text) cells.
$ws.Cells.Item(2, 4).Value = '63.610.17'
$ws.Cells.Item(2, 5).Value = '  -1.48%  '
$ws.Cells.Item(3, 4).Value = '3.039.18'
$ws.Cells.Item(3, 5).Value = '  -2.17%  '
$ws.Cells.Item(4, 5).Value = '  +0.05%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '557.49'
$ws.Cells.Item(5, 5).Value = '  -0.15%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '141.66'
$ws.Cells.Item(6, 5).Value = '  -1.57%  '
$ws.Cells.Item(7, 5).Value = '  +0.01%  '
$ws.Cells.Item(8, 4).Value = '3.033.27'
$ws.Cells.Item(8, 5).Value = '  -2.17%  '
$ws.Cells.Item(9, 5).Value = '  +2.93%  '
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '0.152'
$ws.Cells.Item(10, 5).Value = '  +0.08%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '6.28'
$ws.Cells.Item(11, 5).Value = '  -11.38%  '
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.487'
$ws.Cells.Item(12, 5).Value = '  +5.42%  '
$ws.Cells.Item(13, 5).Value = '  +0.36%  '
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '35.49'
$ws.Cells.Item(14, 5).Value = '  +0.25%  '
$ws.Cells.Item(15, 4).Value = '3.539.53'
$ws.Cells.Item(15, 5).Value = '  -1.80%  '
$ws.Cells.Item(16, 4).Value = '63.701.53'
$ws.Cells.Item(16, 5).Value = '  -1.41%  '
$ws.Cells.Item(17, 4).Value = '3.049.17'
$ws.Cells.Item(17, 5).Value = '  -1.97%  '
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '0.109'
$ws.Cells.Item(18, 5).Value = '  +0.51%  '
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '6.76'
$ws.Cells.Item(19, 5).Value = '  -0.40%  '
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '473.84'
$ws.Cells.Item(20, 5).Value = '  -2.04%  '
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '13.99'
$ws.Cells.Item(21, 5).Value = '  +1.34%  '
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '14.63'
$ws.Cells.Item(22, 5).Value = '  +10.21%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '0.679'
$ws.Cells.Item(23, 5).Value = '  +0.59%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '7.49'
$ws.Cells.Item(24, 5).Value = '  -1.89%  '
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '82.52'
$ws.Cells.Item(25, 5).Value = '  +1.82%  '
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '0.999'
$ws.Cells.Item(26, 5).Value = '  -0.07%  '
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '2.78'
$ws.Cells.Item(27, 5).Value = '  -1.20%  '
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '8.09'
$ws.Cells.Item(28, 5).Value = '  +0.34%  '
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '2.02'
$ws.Cells.Item(29, 5).Value = '  -2.52%  '
$ws.Cells.Item(30, 5).Value = '  +0.01%  '
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '26.09'
$ws.Cells.Item(31, 5).Value = '  -0.27%  '
$ws.Cells.Item(32, 5).Value = '  -1.84%  '
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '2.43'
$ws.Cells.Item(33, 5).Value = '  -0.81%  '
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '5.73'
$ws.Cells.Item(34, 5).Value = '  -0.70%  '
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '6.19'
$ws.Cells.Item(35, 5).Value = '  +0.90%  '
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '54.53'
$ws.Cells.Item(36, 5).Value = '  -1.34%  '
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '0.0408'
$ws.Cells.Item(37, 5).Value = '  -0.26%  '
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '443.34'
$ws.Cells.Item(38, 5).Value = '  -4.51%  '
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.0809'
$ws.Cells.Item(39, 5).Value = '  -2.26%  '
$ws.Cells.Item(40, 2).Value = 'Maker'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(40, 4).Value = '3.015.22'
$ws.Cells.Item(40, 5).Value = '  -0.58%  '
$ws.Cells.Item(41, 2).Value = 'dogwifhat'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '2.78'
$ws.Cells.Item(41, 5).Value = '  +2.75%  '
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '0.117'
$ws.Cells.Item(42, 5).Value = '  +0.20%  '
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '8.24'
$ws.Cells.Item(43, 5).Value = '  -0.59%  '
$ws.Cells.Item(44, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '28.07'
$ws.Cells.Item(44, 5).Value = '  -0.59%  '
$ws.Cells.Item(45, 2).Value = 'TheGraph'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '0.267'
$ws.Cells.Item(45, 5).Value = '  +2.27%  '
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '2.25'
$ws.Cells.Item(46, 5).Value = '  +6.87%  '
$ws.Cells.Item(47, 5).Value = '  -0.03%  '
$ws.Cells.Item(48, 5).Value = '  +0.32%  '
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '117.88'
$ws.Cells.Item(49, 5).Value = '  -0.93%  '
$ws.Cells.Item(50, 4).Value = '0.0₃0510'
$ws.Cells.Item(50, 5).Value = '  -1.40%  '
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '2.07'
$ws.Cells.Item(51, 5).Value = '  -0.10%  '
